$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Rewrite the first body paragraph (paragraph 2) text.
#    Old: "Although some features we have not implemented ... wasted. Therefore, ... client."
#    New: "Although there are some features ... project. Therefore, ... client."
# ---------------------------------------------------------------------------
$old1 = "Although some features we have not implemented such as the login system and the automated status chang"
$new1 = "Although there are some features that we have not implemented such as the login system and the automated status chang"
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

$old2 = "ing mechanism, which we have not done well in this sprint, due to the temporary difficulties we met, we have implemented not only the sprint backlog but also part of the project view and role determination in the project. In this part, we did quite well to change our schedule so that the effort hours were not wasted. Therefore, in next sprint, we are going to finish all the unfinished features and deliver the more complete product to the client."
$new2 = "ing mechanism, we have implemented not only the sprint backlog but also part of the project view and role determination in the project. Therefore, in next sprint, we are going to finish all the unfinished features and deliver the more complete product to the client."
$d.Content.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Temporarily remove the "_GoBack" bookmark; it will be re-created at the
#    true end of the document once all new content has been appended.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 3. Build the paragraph structure first (every new paragraph inherits the
#    plain, non-bold formatting of paragraph 2), THEN go back and apply
#    bold+underline formatting only to the three heading paragraphs. Doing
#    the formatting last prevents it from leaking into paragraphs created
#    afterwards via InsertParagraphAfter.
# ---------------------------------------------------------------------------
$tail = $d.Paragraphs(2).Range
$tail.InsertParagraphAfter() | Out-Null          # 3: blank
$tail = $d.Paragraphs(3).Range
$tail.InsertParagraphAfter() | Out-Null          # 4: heading "What did we do well?"
$tail = $d.Paragraphs(4).Range
$tail.InsertParagraphAfter() | Out-Null          # 5: body
$tail = $d.Paragraphs(5).Range
$tail.InsertParagraphAfter() | Out-Null          # 6: blank
$tail = $d.Paragraphs(6).Range
$tail.InsertParagraphAfter() | Out-Null          # 7: heading "What should we have done better?"
$tail = $d.Paragraphs(7).Range
$tail.InsertParagraphAfter() | Out-Null          # 8: body
$tail = $d.Paragraphs(8).Range
$tail.InsertParagraphAfter() | Out-Null          # 9: blank
$tail = $d.Paragraphs(9).Range
$tail.InsertParagraphAfter() | Out-Null          # 10: heading "Last Retrospective Outcome Result"
$tail = $d.Paragraphs(10).Range
$tail.InsertParagraphAfter() | Out-Null          # 11: body (final)

# Fill in the text for each new paragraph.
$d.Paragraphs(4).Range.Text = "What did we do well?"
$d.Paragraphs(5).Range.Text = "We did quite well about changingthe schedule during development so that the effort hours were not wasted."
$d.Paragraphs(7).Range.Text = "What should we have done better?"
$d.Paragraphs(8).Range.Text = "The time management in this sprint is still not satisfying as we are not really familiar with Django."
$d.Paragraphs(10).Range.Text = "Last Retrospective Outcome Result"
$d.Paragraphs(11).Range.Text = "Fixed as we now have a clearer picture on the final product."

# Apply bold + underline to the three heading paragraphs only.
$d.Paragraphs(4).Range.Font.Bold = 1
$d.Paragraphs(4).Range.Font.Underline = 1
$d.Paragraphs(7).Range.Font.Bold = 1
$d.Paragraphs(7).Range.Font.Underline = 1
$d.Paragraphs(10).Range.Font.Bold = 1
$d.Paragraphs(10).Range.Font.Underline = 1

# ---------------------------------------------------------------------------
# 4. Re-create the "_GoBack" bookmark collapsed at the very end of the
#    document (right after the final run of the last paragraph).
# ---------------------------------------------------------------------------
$endRange = $d.Paragraphs(11).Range.Duplicate
$endRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $endRange) | Out-Null

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
